# Updates "想去人数" (F) / "最低票价" (G) counters across all four sheets and
# rewrites three rows (31-33) on 全部类型 so it catches up with events that
# 演出 already reflects (corresponds to "Update gh-pages to output generated
# at 456a3b4").

$wb = $excel.ActiveWorkbook

# Helper: Excel's COM .Value setter auto-parses plain "YYYY-MM-DD" strings
# into date serials (changing both the stored type and adding a number
# format / style). The source file keeps these as literal text with the
# default style, so force a Text format for the assignment, then drop the
# style back to Normal/General so no stray style index is left behind.
function Set-TextDate($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 809
$ws.Cells.Item(5, 6).Value = 863
$ws.Cells.Item(6, 6).Value = 676
$ws.Cells.Item(9, 6).Value = 837
$ws.Cells.Item(10, 6).Value = 707
$ws.Cells.Item(13, 6).Value = 373
$ws.Cells.Item(14, 6).Value = 727
$ws.Cells.Item(15, 6).Value = 977
$ws.Cells.Item(16, 6).Value = 10210
$ws.Cells.Item(17, 6).Value = 638
$ws.Cells.Item(18, 6).Value = 51
$ws.Cells.Item(21, 6).Value = 47
$ws.Cells.Item(22, 6).Value = 278
$ws.Cells.Item(23, 6).Value = 1780
$ws.Cells.Item(25, 6).Value = 292
$ws.Cells.Item(26, 6).Value = 492
$ws.Cells.Item(27, 6).Value = 188
$ws.Cells.Item(29, 6).Value = 283
$ws.Cells.Item(30, 6).Value = 195
$ws.Cells.Item(32, 6).Value = 76
$ws.Cells.Item(33, 6).Value = 102
$ws.Cells.Item(35, 6).Value = 181
$ws.Cells.Item(36, 6).Value = 201
$ws.Cells.Item(37, 6).Value = 185
$ws.Cells.Item(38, 6).Value = 47
$ws.Cells.Item(39, 6).Value = 96

# --- sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 105
$ws.Cells.Item(7, 6).Value = 139
$ws.Cells.Item(16, 6).Value = 314

# --- sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 826

# --- sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 826
$ws.Cells.Item(6, 6).Value = 809
$ws.Cells.Item(8, 6).Value = 863
$ws.Cells.Item(9, 6).Value = 676
$ws.Cells.Item(12, 6).Value = 105
$ws.Cells.Item(13, 6).Value = 139
$ws.Cells.Item(14, 6).Value = 837
$ws.Cells.Item(15, 6).Value = 707
$ws.Cells.Item(18, 6).Value = 977
$ws.Cells.Item(19, 6).Value = 10214
$ws.Cells.Item(21, 6).Value = 638
$ws.Cells.Item(22, 6).Value = 51
$ws.Cells.Item(24, 6).Value = 278
$ws.Cells.Item(25, 6).Value = 1780
$ws.Cells.Item(26, 6).Value = 492
$ws.Cells.Item(27, 6).Value = 188

# Row 31: was "广州·星宇宙动漫嘉年华" -> becomes the wanuka gig
# (shifts the whole row's content the way 演出 rows 14-16 already show it).
Set-TextDate $ws.Cells.Item(31, 2) "2024-09-01"
$ws.Cells.Item(31, 3).Value = " 广州·日本次世代神秘创作歌手 和ぬか(wanuka) 演出"
$ws.Cells.Item(31, 4).Value = "南洲路158号2F SD Livehouse"
$ws.Cells.Item(31, 5).Value = "2024.09.01 20:00-09.01 22:00"
$ws.Cells.Item(31, 6).Value = 5
$ws.Cells.Item(31, 7).Value = 260
$ws.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90116"
$ws.Cells.Item(31, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/6VG7pF6k1722407692893.jpeg"

# Row 32: was the wanuka gig -> becomes the 音阅派国漫演唱会 concert.
$ws.Cells.Item(32, 3).Value = "广州·音阅派国漫演唱会-《狐妖小红娘》《一人之下》领衔国漫原声音乐现场"
$ws.Cells.Item(32, 4).Value = "东风中路259号 广州中山纪念堂"
$ws.Cells.Item(32, 5).Value = "2024.09.01 19:30-09.01 21:00"
$ws.Cells.Item(32, 6).Value = 63
$ws.Cells.Item(32, 7).Value = 180
$ws.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89794"
$ws.Cells.Item(32, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/ehol1aeZ1721810539643.jpeg"

# Row 33: was the 音阅派 concert -> becomes 神山羊2024巡演 ENCOUNTER.
Set-TextDate $ws.Cells.Item(33, 2) "2024-09-08"
$ws.Cells.Item(33, 3).Value = "广州·神山羊2024巡演ENCOUNTER"
$ws.Cells.Item(33, 4).Value = "流花路117号流花展贸中心5号馆 广州大麦66live house"
$ws.Cells.Item(33, 5).Value = "2024.09.08 19:00-09.08 20:30"
$ws.Cells.Item(33, 6).Value = 314
$ws.Cells.Item(33, 7).Value = 380
$ws.Cells.Item(33, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89835"
$ws.Cells.Item(33, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/vpWr6GKQ1721877449091.jpeg"

$ws.Cells.Item(34, 6).Value = 283
$ws.Cells.Item(36, 6).Value = 195
$ws.Cells.Item(38, 6).Value = 76
$ws.Cells.Item(39, 6).Value = 102
$ws.Cells.Item(42, 6).Value = 181
$ws.Cells.Item(45, 6).Value = 201
$ws.Cells.Item(46, 6).Value = 185
